$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 17.93638541084482
$ws.Cells.Item(2, 3).Value = 6.267222904724992
$ws.Cells.Item(2, 4).Value = 3.776792924455576
$ws.Cells.Item(2, 5).Value = 10.49183963931405
$ws.Cells.Item(2, 6).Value = 56.73234970993191
$ws.Cells.Item(2, 9).Value = 41.13147533566444
$ws.Cells.Item(2, 10).Value = 10.31861299239345
$ws.Cells.Item(2, 11).Value = 16.29887471242854
$ws.Cells.Item(2, 12).Value = 11.88618270867793

$ws.Cells.Item(3, 2).Value = 17.89379833172608
$ws.Cells.Item(3, 3).Value = 6.147647537142801
$ws.Cells.Item(3, 4).Value = 3.795759176705817
$ws.Cells.Item(3, 5).Value = 10.51373394722378
$ws.Cells.Item(3, 6).Value = 56.54239607837059
$ws.Cells.Item(3, 9).Value = 41.03273893390486
$ws.Cells.Item(3, 10).Value = 10.33078947304548
$ws.Cells.Item(3, 11).Value = 16.27115215239635
$ws.Cells.Item(3, 12).Value = 11.90546665009949

$ws.Cells.Item(4, 2).Value = 17.87280340412051
$ws.Cells.Item(4, 3).Value = 6.075128913113737
$ws.Cells.Item(4, 4).Value = 3.808166578393161
$ws.Cells.Item(4, 5).Value = 10.52828887710655
$ws.Cells.Item(4, 6).Value = 56.43286028250733
$ws.Cells.Item(4, 9).Value = 40.97632699897954
$ws.Cells.Item(4, 10).Value = 10.33886500076158
$ws.Cells.Item(4, 11).Value = 16.2584155173799
$ws.Cells.Item(4, 12).Value = 11.91909959623436

$ws.Cells.Item(5, 2).Value = 17.8655519223222
$ws.Cells.Item(5, 3).Value = 6.045848442110091
$ws.Cells.Item(5, 4).Value = 3.81341485492007
$ws.Cells.Item(5, 5).Value = 10.53450023911083
$ws.Cells.Item(5, 6).Value = 56.39002711416671
$ws.Cells.Item(5, 9).Value = 40.95440529428512
$ws.Cells.Item(5, 10).Value = 10.3423068015808
$ws.Cells.Item(5, 11).Value = 16.25430763588231
$ws.Cells.Item(5, 12).Value = 11.92510608772493

$ws.Cells.Item(6, 2).Value = 17.86442679583351
$ws.Cells.Item(6, 3).Value = 6.041004148198087
$ws.Cells.Item(6, 4).Value = 3.814297949058034
$ws.Cells.Item(6, 5).Value = 10.53554856526478
$ws.Cells.Item(6, 6).Value = 56.38302415199504
$ws.Cells.Item(6, 9).Value = 40.95082983117467
$ws.Cells.Item(6, 10).Value = 10.34288743639976
$ws.Cells.Item(6, 11).Value = 16.25369102236763
$ws.Cells.Item(6, 12).Value = 11.92613070279154

$ws.Cells.Item(7, 2).Value = 17.87270031805633
$ws.Cells.Item(7, 3).Value = 6.074732868051564
$ws.Cells.Item(7, 4).Value = 3.808236579714547
$ws.Cells.Item(7, 5).Value = 10.52837151081204
$ws.Cells.Item(7, 6).Value = 56.43227529315909
$ws.Cells.Item(7, 9).Value = 40.97602702754244
$ws.Cells.Item(7, 10).Value = 10.33891080646598
$ws.Cells.Item(7, 11).Value = 16.25835572878965
$ws.Cells.Item(7, 12).Value = 11.91917877574259

$ws.Cells.Item(8, 2).Value = 17.92063707169216
$ws.Cells.Item(8, 3).Value = 6.225833167128649
$ws.Cells.Item(8, 4).Value = 3.783174647639991
$ws.Cells.Item(8, 5).Value = 10.49915838940633
$ws.Cells.Item(8, 6).Value = 56.66539277675376
$ws.Cells.Item(8, 9).Value = 41.09656130428834
$ws.Cells.Item(8, 10).Value = 10.3226872784836
$ws.Cells.Item(8, 11).Value = 16.28842981370651
$ws.Cells.Item(8, 12).Value = 11.89245991162408

$ws.Cells.Item(9, 2).Value = 18.05510837660585
$ws.Cells.Item(9, 3).Value = 6.527333364683534
$ws.Cells.Item(9, 4).Value = 3.740050810041211
$ws.Cells.Item(9, 5).Value = 10.45066811854126
$ws.Cells.Item(9, 6).Value = 57.17793648595561
$ws.Cells.Item(9, 9).Value = 41.36600640648416
$ws.Cells.Item(9, 10).Value = 10.29561326735543
$ws.Cells.Item(9, 11).Value = 16.38114755841291
$ws.Cells.Item(9, 12).Value = 11.8542769554583

$ws.Cells.Item(10, 2).Value = 18.17791207231296
$ws.Cells.Item(10, 3).Value = 6.749443734800346
$ws.Cells.Item(10, 4).Value = 3.712006463408853
$ws.Cells.Item(10, 5).Value = 10.42037164457624
$ws.Cells.Item(10, 6).Value = 57.58696627340885
$ws.Cells.Item(10, 9).Value = 41.58363032356525
$ws.Cells.Item(10, 10).Value = 10.2785935422438
$ws.Cells.Item(10, 11).Value = 16.46942855084456
$ws.Cells.Item(10, 12).Value = 11.83487293966059

$ws.Cells.Item(11, 2).Value = 18.23882206970751
$ws.Cells.Item(11, 3).Value = 6.850098566240064
$ws.Cells.Item(11, 4).Value = 3.700031606282018
$ws.Cells.Item(11, 5).Value = 10.40773914105901
$ws.Cells.Item(11, 6).Value = 57.77978099216224
$ws.Cells.Item(11, 9).Value = 41.68678013741296
$ws.Cells.Item(11, 10).Value = 10.27147050443699
$ws.Cells.Item(11, 11).Value = 16.51386016131464
$ws.Cells.Item(11, 12).Value = 11.82791925574072

$ws.Cells.Item(12, 2).Value = 18.26259641209464
$ws.Cells.Item(12, 3).Value = 6.888117779970691
$ws.Cells.Item(12, 4).Value = 3.695609056831714
$ws.Cells.Item(12, 5).Value = 10.40312027599728
$ws.Cells.Item(12, 6).Value = 57.85373559035433
$ws.Cells.Item(12, 9).Value = 41.72642542281078
$ws.Cells.Item(12, 10).Value = 10.26886194863951
$ws.Cells.Item(12, 11).Value = 16.53128899963657
$ws.Cells.Item(12, 12).Value = 11.82555497256775

$ws.Cells.Item(13, 2).Value = 18.2574449233359
$ws.Cells.Item(13, 3).Value = 6.879934593807578
$ws.Cells.Item(13, 4).Value = 3.696556555198409
$ws.Cells.Item(13, 5).Value = 10.40410771026639
$ws.Cells.Item(13, 6).Value = 57.83776683860074
$ws.Cells.Item(13, 9).Value = 41.71786130067465
$ws.Cells.Item(13, 10).Value = 10.26941980345435
$ws.Cells.Item(13, 11).Value = 16.52750872552905
$ws.Cells.Item(13, 12).Value = 11.82605221058178

$ws.Cells.Item(14, 2).Value = 18.24076388788504
$ws.Cells.Item(14, 3).Value = 6.853228594962458
$ws.Cells.Item(14, 4).Value = 3.699665517149938
$ws.Cells.Item(14, 5).Value = 10.40735584453239
$ws.Cells.Item(14, 6).Value = 57.78584662291711
$ws.Cells.Item(14, 9).Value = 41.6900301339306
$ws.Cells.Item(14, 10).Value = 10.27125411918038
$ws.Cells.Item(14, 11).Value = 16.51528201120028
$ws.Cells.Item(14, 12).Value = 11.82771935788276

$ws.Cells.Item(15, 2).Value = 18.2306380990876
$ws.Cells.Item(15, 3).Value = 6.836856629661704
$ws.Cells.Item(15, 4).Value = 3.701584426797028
$ws.Cells.Item(15, 5).Value = 10.40936686595316
$ws.Cells.Item(15, 6).Value = 57.75416547029032
$ws.Cells.Item(15, 9).Value = 41.67305850685468
$ws.Cells.Item(15, 10).Value = 10.27238924472855
$ws.Cells.Item(15, 11).Value = 16.5078710605739
$ws.Cells.Item(15, 12).Value = 11.82877554132679

$ws.Cells.Item(16, 2).Value = 18.17403161112663
$ws.Cells.Item(16, 3).Value = 6.742854393032315
$ws.Cells.Item(16, 4).Value = 3.7128047567002
$ws.Cells.Item(16, 5).Value = 10.42122029606329
$ws.Cells.Item(16, 6).Value = 57.57449856753736
$ws.Cells.Item(16, 9).Value = 41.57697174631046
$ws.Cells.Item(16, 10).Value = 10.27907148842499
$ws.Cells.Item(16, 11).Value = 16.46660990660891
$ws.Cells.Item(16, 12).Value = 11.83536503765181

$ws.Cells.Item(17, 2).Value = 18.14058636894487
$ws.Cells.Item(17, 3).Value = 6.685058600093879
$ws.Cells.Item(17, 4).Value = 3.719888182085686
$ws.Cells.Item(17, 5).Value = 10.42878604770298
$ws.Cells.Item(17, 6).Value = 57.4659853727728
$ws.Cells.Item(17, 9).Value = 41.51908043435193
$ws.Cells.Item(17, 10).Value = 10.2833292563951
$ws.Cells.Item(17, 11).Value = 16.44238434666205
$ws.Cells.Item(17, 12).Value = 11.83988697583987

$ws.Cells.Item(18, 2).Value = 18.12182557264948
$ws.Cells.Item(18, 3).Value = 6.651782167769419
$ws.Cells.Item(18, 4).Value = 3.724036076451825
$ws.Cells.Item(18, 5).Value = 10.43324589615136
$ws.Cells.Item(18, 6).Value = 57.40420822240798
$ws.Cells.Item(18, 9).Value = 41.48617443532436
$ws.Cells.Item(18, 10).Value = 10.28583652159873
$ws.Cells.Item(18, 11).Value = 16.4288532641465
$ws.Cells.Item(18, 12).Value = 11.84266423009667

$ws.Cells.Item(19, 2).Value = 18.1155557311432
$ws.Cells.Item(19, 3).Value = 6.640510822995675
$ws.Cells.Item(19, 4).Value = 3.725453154083107
$ws.Cells.Item(19, 5).Value = 10.43477452799828
$ws.Cells.Item(19, 6).Value = 57.38340179867738
$ws.Cells.Item(19, 9).Value = 41.47510064180392
$ws.Cells.Item(19, 10).Value = 10.28669546160428
$ws.Cells.Item(19, 11).Value = 16.42434136998063
$ws.Cells.Item(19, 12).Value = 11.84363486252535

$ws.Cells.Item(20, 2).Value = 18.1440975217426
$ws.Cells.Item(20, 3).Value = 6.691214867076183
$ws.Cells.Item(20, 4).Value = 3.719126515630215
$ws.Cells.Item(20, 5).Value = 10.42796946301071
$ws.Cells.Item(20, 6).Value = 57.47747109789132
$ws.Cells.Item(20, 9).Value = 41.52520262614298
$ws.Cells.Item(20, 10).Value = 10.28286997689009
$ws.Cells.Item(20, 11).Value = 16.44492158391619
$ws.Cells.Item(20, 12).Value = 11.83938735882692

$ws.Cells.Item(21, 2).Value = 18.24564440688548
$ws.Cells.Item(21, 3).Value = 6.861075734654603
$ws.Cells.Item(21, 4).Value = 3.698749301902847
$ws.Cells.Item(21, 5).Value = 10.40639732047089
$ws.Cells.Item(21, 6).Value = 57.80107158845875
$ws.Cells.Item(21, 9).Value = 41.69818905560673
$ws.Cells.Item(21, 10).Value = 10.27071292878873
$ws.Cells.Item(21, 11).Value = 16.51885700372516
$ws.Cells.Item(21, 12).Value = 11.82722238165701

$ws.Cells.Item(22, 2).Value = 18.31613553321995
$ws.Cells.Item(22, 3).Value = 6.971508969194801
$ws.Cells.Item(22, 4).Value = 3.686084622665606
$ws.Cells.Item(22, 5).Value = 10.39325896328184
$ws.Cells.Item(22, 6).Value = 58.01802642340084
$ws.Cells.Item(22, 9).Value = 41.81464749031713
$ws.Cells.Item(22, 10).Value = 10.26328497015265
$ws.Cells.Item(22, 11).Value = 16.57069040892204
$ws.Cells.Item(22, 12).Value = 11.82083911696198

$ws.Cells.Item(23, 2).Value = 18.27814148866851
$ws.Cells.Item(23, 3).Value = 6.91263491745187
$ws.Cells.Item(23, 4).Value = 3.692784405651961
$ws.Cells.Item(23, 5).Value = 10.40018345386398
$ws.Cells.Item(23, 6).Value = 57.90174396628822
$ws.Cells.Item(23, 9).Value = 41.7521844187638
$ws.Cells.Item(23, 10).Value = 10.2672021617225
$ws.Cells.Item(23, 11).Value = 16.54270835793963
$ws.Cells.Item(23, 12).Value = 11.82410274798777

$ws.Cells.Item(24, 2).Value = 18.14250867268132
$ws.Cells.Item(24, 3).Value = 6.688431767806061
$ws.Cells.Item(24, 4).Value = 3.719470629714211
$ws.Cells.Item(24, 5).Value = 10.42833829766901
$ws.Cells.Item(24, 6).Value = 57.47227650372945
$ws.Cells.Item(24, 9).Value = 41.52243360954176
$ws.Cells.Item(24, 10).Value = 10.28307743205615
$ws.Cells.Item(24, 11).Value = 16.44377326305317
$ws.Cells.Item(24, 12).Value = 11.83961268266021

$ws.Cells.Item(25, 2).Value = 18.01446088613152
$ws.Cells.Item(25, 3).Value = 6.445483101346213
$ws.Cells.Item(25, 4).Value = 3.75107560713835
$ws.Cells.Item(25, 5).Value = 10.46284772993057
$ws.Cells.Item(25, 6).Value = 57.03347723196411
$ws.Cells.Item(25, 9).Value = 41.28962188535014
$ws.Cells.Item(25, 10).Value = 10.30243191578426
$ws.Cells.Item(25, 11).Value = 16.35249087744309
$ws.Cells.Item(25, 12).Value = 11.86308625632313
